$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 541 ("Fruta / hortaliza, semanal"),
# which pushes the former rows 541-649 down to 542-650.
$ws.Range("A541").EntireRow.Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A541").Value = 3
$ws.Range("B541").Value = "Femacal de La Calera"
$ws.Range("C541").Value = "Coquimbo"
$ws.Range("D541").Value = 45209
$ws.Range("E541").Value = 5
$ws.Range("F541").Value = 100112043
$ws.Range("G541").Value = "Pepino ensalada"
$ws.Range("H541").Value = "Sin especificar"
$ws.Range("I541").Value = "Primera"
$ws.Range("J541").Value = 90
$ws.Range("K541").Value = 13500
$ws.Range("L541").Value = 14000
$ws.Range("M541").Value = 13722
$ws.Range("N541").Value = "$/caja 60 unidades"
$ws.Range("O541").Value = "Región de Arica y Parinacota"
$ws.Range("P541").Value = 229
$ws.Range("Q541").Value = 60
$ws.Range("R541").Value = "Hortaliza"
